$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.106.05"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.827.08"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'313.01"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "'0.4577"
$ws.Range("E7").Value = "  +7.06%  "
$ws.Range("D8").Value = "'0.3737"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("D9").Value = "'0.07317"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").Value = "'0.8632"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.823.41"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'6.718"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "'92.94"
$ws.Range("E15").Value = "  +5.06%  "
$ws.Range("D16").Value = "'0.07089"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "'0.000008861"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "27.125.34"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'5.204"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D24").Value = "'2.004"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'152.05"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "'2.227"
$ws.Range("E26").Value = "  +5.37%  "
$ws.Range("D27").Value = "'18.49"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("D28").Value = "'5.293"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").Value = "'117.61"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("D30").Value = "'0.08897"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").Value = "'0.7622"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").Value = "'2.975"
$ws.Range("E33").Value = "  +5.50%  "
$ws.Range("D34").Value = "'4.482"
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").Value = "'1.106"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").Value = "'0.05302"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").Value = "'0.5375"
$ws.Range("E39").Value = "  +6.95%  "
$ws.Range("D40").Value = "'7.196"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "'2.879"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").Value = "'0.1720"
$ws.Range("D43").Value = "'0.5223"
$ws.Range("E43").Value = "  +11.18%  "
$ws.Range("D44").Value = "'8.636"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "'10.70"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").Value = "'1.987"
$ws.Range("E46").Value = "  +10.32%  "
$ws.Range("D47").Value = "'106.10"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "'1.686"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "'0.06423"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "'0.9244"
$ws.Range("E51").Value = "  +1.30%  "
